# refatoracao e ajustes do tdd
# Update the "Usuario" test value on the Cadastro sheet and move the
# saved cell selection back onto that row (A2) instead of the old A9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cadastro")
$ws.Activate()

# Bump the username used by the TDD suite.
$ws.Range("A2").Value = "pradov1057"

# Leave the cursor/selection on A2 (was A9 before this refactor).
$ws.Range("A2").Select()
